# Tijdschrijfformulier.xlsx update
# - Add a new "KBS b les" time entry (120 min, 2023-11-15, Opmerkingen "Les")
#   to Marvin, Demi, Lucas and Luuk sheets (row 14).
# - On the Jochem sheet: correct the activity on row 15 from "KBS a les" to
#   "KBS b les", and add a new row 16 for "Functionaliteit afrekenpagina"
#   (60 min, 2023-11-15).
# - Make "Marvin" the active sheet/tab again.

$wb = $excel.ActiveWorkbook

function Add-TimeRow {
    param(
        $ws,
        [int]$row,
        [string]$activiteit,
        [int]$datumSerial,
        [int]$minuten,
        [string]$opmerkingen
    )

    # Column A - activity (shared string)
    $ws.Cells.Item($row, 1).Value = $activiteit

    # Column B - date; copy the number format from an existing date cell
    # (B10) so it reuses the workbook's existing date style instead of
    # creating a brand new number format.
    $ws.Range("B10").Copy()
    $ws.Cells.Item($row, 2).PasteSpecial(-4122)
    $ws.Cells.Item($row, 2).Value = $datumSerial

    # Column C - minutes spent
    $ws.Cells.Item($row, 3).Value = $minuten

    # Column D - Opmerkingen (optional)
    if ($opmerkingen) {
        $ws.Cells.Item($row, 4).Value = $opmerkingen
    }
}

# --- Marvin (sheet1) ---------------------------------------------------
$wsMarvin = $wb.Worksheets.Item("Marvin")
Add-TimeRow $wsMarvin 14 "KBS b les" 45245 120 "Les"

# --- Demi (sheet2) -------------------------------------------------------
$wsDemi = $wb.Worksheets.Item("Demi")
Add-TimeRow $wsDemi 14 "KBS b les" 45245 120 "Les"

# --- Lucas (sheet3) ------------------------------------------------------
$wsLucas = $wb.Worksheets.Item("Lucas")
Add-TimeRow $wsLucas 14 "KBS b les" 45245 120 "Les"

# --- Luuk (sheet4) -------------------------------------------------------
$wsLuuk = $wb.Worksheets.Item("Luuk")
Add-TimeRow $wsLuuk 14 "KBS b les" 45245 120 "Les"

# --- Jochem (sheet5) -------------------------------------------------------
$wsJochem = $wb.Worksheets.Item("Jochem")
# Row 15 was mislabeled "KBS a les" - fix to "KBS b les".
$wsJochem.Cells.Item(15, 1).Value = "KBS b les"
# New row 16: Functionaliteit afrekenpagina, 60 minutes, no Opmerkingen.
Add-TimeRow $wsJochem 16 "Functionaliteit afrekenpagina" 45245 60 $null

# --- Selection bookkeeping (mirrors what Excel records when the edits are
#     made interactively - the newly entered row gets selected on each
#     sheet, and the last sheet touched/activated becomes the active tab).
[void]$wsDemi.Range("A14:D14").Select()
[void]$wsLucas.Range("A14:D14").Select()
[void]$wsLuuk.Range("A14:D14").Select()
[void]$wsJochem.Range("A15:D15").Select()
[void]$wsMarvin.Range("A14:D14").Select()
$wsMarvin.Activate()
